# Apply the "notebook can generate ppt now" edit:
#  - Slide 1: retitle + collapse the 4-paragraph summary into a single line.
#  - Slides 2-10: retitle + replace the bodies with new "-" prefixed bullets.
#  - Slides 11-13: removed entirely (deck shrinks from 13 to 10 slides).

$p = $ppt.ActivePresentation

# --- Remove the trailing three slides (delete from the end so indices stay valid) ---
$p.Slides.Item(13).Delete()
$p.Slides.Item(12).Delete()
$p.Slides.Item(11).Delete()

# --- Slide 1: Title ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Title"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Learning with Experts for Fine-grained Category Discovery"

# --- Slide 2: Introduction ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "- Addressing the problem of generalized category discovery (GCD)`r- Introducing Expert-Contrastive Learning (XCon) method for mining useful information from images`r- Utilizing k-means clustering and contrastive learning on sub-datasets for learning discriminative features"

# --- Slide 3: Problem Statement ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Problem Statement"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "- GCD involves discovering categories within unlabeled data by leveraging information`r- Existing methods tend to cluster based on class-irrelevant cues, leading to suboptimal results`r- Need for fine-grained category discovery with a focus on relevant concepts"

# --- Slide 4: XCon Methodology ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "XCon Methodology"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "- Partitioning data into k expert sub-datasets using k-means clustering on self-supervised representations`r- Each sub-dataset acts as an expert dataset to eliminate negative influences of class-irrelevant cues`r- Learning discriminative features for fine-grained category discovery"

# --- Slide 5: Novel Category Discovery ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Novel Category Discovery"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "- Objective: Discover new object categories by transferring knowledge learned from seen classes`r- XCon partitions data into k sub-datasets for learning discriminative representations`r- Setting a new state-of-the-art performance on tested category discovery benchmarks"

# --- Slide 6: Previous Work ---
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Previous Work"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "- Contrastive learning in NCD problem by NCL`r- Using k-means grouping on self-supervised features for informative contrastive pairs`r- Focus on Generalized Category Discovery for effective representation learning"

# --- Slide 7: Experimental Results (Generic Datasets) ---
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Experimental Results (Generic Datasets)"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "- Results on CIFAR10, CIFAR100, Stanford Cars, etc.`r- XCon outperforms baseline methods, showing consistent improvement`r- Evaluation metric: Clustering accuracy (ACC) on unlabeled dataset"

# --- Slide 8: Experimental Results (Fine-grained Datasets) ---
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Experimental Results (Fine-grained Datasets)"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "- Performance on fine-grained image classification benchmarks`r- Improved ACC on CUB-200 and Stanford Cars with XCon`r- Analysis of weight parameter α for fine-grained loss"

# --- Slide 9: Qualitative Analysis ---
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Qualitative Analysis"
$s9.Shapes.Item(2).TextFrame.TextRange.Text = "- Visualization of feature space with XCon compared to DINO`r- Clear boundaries between groups corresponding to different categories`r- Demonstrating discriminative features learned by XCon"

# --- Slide 10: Conclusion ---
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "- XCon addresses GCD with self-supervised representation`r- Achieving improved performance in fine-grained category discovery`r- Validation of method effectiveness through experiments and comparisons"
